$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (Question) and C (ConditionType), rows 3-21
# (row 2 is unchanged)
$data = @{
    3  = @(7, 2)
    4  = @(4, 2)
    5  = @(3, 2)
    6  = @(10, 1)
    7  = @(9, 2)
    8  = @(8, 1)
    9  = @(5, 1)
    10 = @(1, 1)
    11 = @(2, 2)
    12 = @(3, 2)
    13 = @(5, 1)
    14 = @(1, 1)
    15 = @(6, 1)
    16 = @(10, 2)
    17 = @(2, 1)
    18 = @(7, 1)
    19 = @(4, 1)
    20 = @(8, 2)
    21 = @(9, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}

# Update the selection to match the edited workbook state
$ws.Range("A1:C11").Select()
